# Generate Report for Handback
#
# The localization status report is regenerated: both the zh-cn and de-de
# targets have now been handed back and are in sync with en-US, so:
#   - the "Status" column moves from "Ready for handoff" to
#     "Handed back: in sync with en-US" (on the Overview sheet and on each
#     per-language sheet)
#   - the "Latest Handback DateTime" is refreshed to the new handback time
#   - the stale "Error Detail" (out-of-date handback file warning) is cleared
#     now that the handback is current

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ---- Overview sheet -------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus

# Column widths auto-grow to fit the longer status text.
$overview.Columns.Item(5).ColumnWidth = 29.166666666666668
$overview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---- zh-cn sheet ------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("K2").Value = "2016-08-13 15:00:43"
$zhcn.Range("P2").Value = ""

$zhcn.Columns.Item(3).ColumnWidth = 29.166666666666668
$zhcn.Columns.Item(16).ColumnWidth = 12.833333333333332

# ---- de-de sheet ------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Range("K2").Value = "2016-08-13 15:00:52"
$dede.Range("P2").Value = ""

$dede.Columns.Item(3).ColumnWidth = 29.166666666666668
$dede.Columns.Item(16).ColumnWidth = 12.833333333333332
